$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029768215540761
$ws.Range("D2").Value = 1.034140409772263
$ws.Range("E2").Value = 1.039107914137565
$ws.Range("F2").Value = 1.049100804991565
$ws.Range("I2").Value = 1.034551686873269
$ws.Range("J2").Value = 1.034913164017445
$ws.Range("K2").Value = 1.036940739213756
$ws.Range("L2").Value = 1.041894040855497
$ws.Range("M2").Value = 1.051858802697866
$ws.Range("N2").Value = 1.015640517648769

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030562604720285
$ws.Range("D3").Value = 1.034729235998363
$ws.Range("E3").Value = 1.039887838461232
$ws.Range("F3").Value = 1.050130206903502
$ws.Range("I3").Value = 1.034704572261864
$ws.Range("J3").Value = 1.035349517361032
$ws.Range("K3").Value = 1.037339278951794
$ws.Range("L3").Value = 1.042484195792342
$ws.Range("M3").Value = 1.052699821718335
$ws.Range("N3").Value = 1.015786088068981

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03107719075133
$ws.Range("D4").Value = 1.035110658128329
$ws.Range("E4").Value = 1.040393457374549
$ws.Range("F4").Value = 1.050797798781104
$ws.Range("I4").Value = 1.034802553679574
$ws.Range("J4").Value = 1.035631757134885
$ws.Range("K4").Value = 1.03759687396313
$ws.Range("L4").Value = 1.042866362008211
$ws.Range("M4").Value = 1.053244908488761
$ws.Range("N4").Value = 1.015880215133839

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031293655879763
$ws.Range("D5").Value = 1.035271104928151
$ws.Range("E5").Value = 1.04060624665822
$ws.Range("F5").Value = 1.051078811844094
$ws.Range("I5").Value = 1.034843518014146
$ws.Range("J5").Value = 1.035750383222436
$ws.Range("K5").Value = 1.037705097018988
$ws.Range("L5").Value = 1.043027094324917
$ws.Range("M5").Value = 1.053474274798011
$ws.Range("N5").Value = 1.015919769746304

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031330009088049
$ws.Range("D6").Value = 1.035298050285533
$ws.Range("E6").Value = 1.040641988182392
$ws.Range("F6").Value = 1.051126016085846
$ws.Range("I6").Value = 1.034850382772423
$ws.Range("J6").Value = 1.035770299416969
$ws.Range("K6").Value = 1.037723264012752
$ws.Range("L6").Value = 1.043054086032762
$ws.Range("M6").Value = 1.053512798788049
$ws.Range("N6").Value = 1.015926410165246

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.0310800826478
$ws.Range("D7").Value = 1.035112801647619
$ws.Range("E7").Value = 1.040396299785197
$ws.Range("F7").Value = 1.050801552288999
$ws.Range("I7").Value = 1.034803101940006
$ws.Range("J7").Value = 1.035633342332053
$ws.Range("K7").Value = 1.037598320321031
$ws.Range("L7").Value = 1.042868509449342
$ws.Range("M7").Value = 1.053247972462171
$ws.Range("N7").Value = 1.01588074372951

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030036565494717
$ws.Range("D8").Value = 1.034339319959183
$ws.Range("E8").Value = 1.03937129429718
$ws.Range("F8").Value = 1.049448384597369
$ws.Range("I8").Value = 1.034603550610861
$ws.Range("J8").Value = 1.035060653746207
$ws.Range("K8").Value = 1.037075486247788
$ws.Range("L8").Value = 1.042093424275998
$ws.Range("M8").Value = 1.052142843557688
$ws.Range("N8").Value = 1.015689727412431

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028202140733099
$ws.Range("D9").Value = 1.032979580252644
$ws.Range("E9").Value = 1.037572499021137
$ws.Range("F9").Value = 1.047075492371678
$ws.Range("I9").Value = 1.034244702462195
$ws.Range("J9").Value = 1.034050712659707
$ws.Range("K9").Value = 1.036152042611804
$ws.Range("L9").Value = 1.040729955623369
$ws.Range("M9").Value = 1.050202358773773
$ws.Range("N9").Value = 1.015352640186431

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026982237988776
$ws.Range("D10").Value = 1.032075365358866
$ws.Range("E10").Value = 1.036378376420158
$ws.Range("F10").Value = 1.045501437308806
$ws.Range("I10").Value = 1.034000660220319
$ws.Range("J10").Value = 1.033376957114883
$ws.Range("K10").Value = 1.035535045044664
$ws.Range("L10").Value = 1.039822623644705
$ws.Range("M10").Value = 1.048913428347274
$ws.Range("N10").Value = 1.015127610167481

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026454750066672
$ws.Range("D11").Value = 1.031684392695804
$ws.Range("E11").Value = 1.035862533230381
$ws.Range("F11").Value = 1.044821742356645
$ws.Range("I11").Value = 1.033893854684054
$ws.Range("J11").Value = 1.033085118847763
$ws.Range("K11").Value = 1.035267569004672
$ws.Range("L11").Value = 1.039430146472739
$ws.Range("M11").Value = 1.048356447403882
$ws.Range("N11").Value = 1.015030102703999

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02625893010772
$ws.Range("D12").Value = 1.031539253711349
$ws.Range("E12").Value = 1.035671111102089
$ws.Range("F12").Value = 1.044569557707585
$ws.Range("I12").Value = 1.033854012695637
$ws.Range("J12").Value = 1.032976703634091
$ws.Range("K12").Value = 1.035168170845543
$ws.Range("L12").Value = 1.039284424967003
$ws.Range("M12").Value = 1.048149731859914
$ws.Range("N12").Value = 1.014993874311898

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026300929034162
$ws.Range("D13").Value = 1.031570382610803
$ws.Range("E13").Value = 1.035712163398765
$ws.Range("F13").Value = 1.044623639277797
$ws.Range("I13").Value = 1.033862566604124
$ws.Range("J13").Value = 1.032999959656714
$ws.Range("K13").Value = 1.035189494134048
$ws.Range("L13").Value = 1.039315679897568
$ws.Range("M13").Value = 1.048194065238647
$ws.Range("N13").Value = 1.015001645860513

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026438561222333
$ws.Range("D14").Value = 1.031672393704108
$ws.Range("E14").Value = 1.035846706433511
$ws.Range("F14").Value = 1.044800890885122
$ws.Range("I14").Value = 1.03389056478999
$ws.Range("J14").Value = 1.033076157480651
$ws.Range("K14").Value = 1.035259353645938
$ws.Range("L14").Value = 1.039418099815147
$ws.Range("M14").Value = 1.048339356715302
$ws.Range("N14").Value = 1.015027108250775

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026523375918567
$ws.Range("D15").Value = 1.031735257511003
$ws.Range("E15").Value = 1.035929627404173
$ws.Range("F15").Value = 1.044910139175852
$ws.Range("I15").Value = 1.033907792934942
$ws.Range("J15").Value = 1.033123103723147
$ws.Range("K15").Value = 1.035302390384765
$ws.Range("L15").Value = 1.039481212361773
$ws.Range("M15").Value = 1.04842889842423
$ws.Range("N15").Value = 1.015042795187011

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02701726131598
$ws.Range("D16").Value = 1.032101324843118
$ws.Range("E16").Value = 1.036412637081803
$ws.Range("F16").Value = 1.045546586203803
$ws.Range("I16").Value = 1.03400772472384
$ws.Range("J16").Value = 1.033396323503011
$ws.Range("K16").Value = 1.035552790079517
$ws.Range("L16").Value = 1.039848679680691
$ws.Range("M16").Value = 1.0489504173709
$ws.Range("N16").Value = 1.015134080021538

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027327261338588
$ws.Range("D17").Value = 1.032331099956726
$ws.Range("E17").Value = 1.036715944090437
$ws.Range("F17").Value = 1.045946317464662
$ws.Range("I17").Value = 1.034070106253299
$ws.Range("J17").Value = 1.03356768174105
$ws.Range("K17").Value = 1.035709776523273
$ws.Range("L17").Value = 1.040079291226392
$ws.Range("M17").Value = 1.049277857235824
$ws.Range("N17").Value = 1.015191322693163

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027508150203183
$ws.Range("D18").Value = 1.032465177646448
$ws.Range("E18").Value = 1.036892975510858
$ws.Range("F18").Value = 1.046179655349293
$ws.Range("I18").Value = 1.034106382920374
$ws.Range("J18").Value = 1.033667622598792
$ws.Range("K18").Value = 1.035801313871902
$ws.Range("L18").Value = 1.040213841885959
$ws.Range("M18").Value = 1.049468956722289
$ws.Range("N18").Value = 1.015224704786805

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027569840640014
$ws.Range("D19").Value = 1.03251090373695
$ws.Range("E19").Value = 1.036953358539874
$ws.Range("F19").Value = 1.046259248256072
$ws.Range("I19").Value = 1.034118733746198
$ws.Range("J19").Value = 1.03370169820055
$ws.Range("K19").Value = 1.035832520584437
$ws.Range("L19").Value = 1.040259726731466
$ws.Range("M19").Value = 1.049534135210063
$ws.Range("N19").Value = 1.015236086087358

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027293993904162
$ws.Range("D20").Value = 1.032306441688761
$ws.Range("E20").Value = 1.036683389945094
$ws.Range("F20").Value = 1.045903411307901
$ws.Range("I20").Value = 1.034063424619011
$ws.Range("J20").Value = 1.033549297598242
$ws.Range("K20").Value = 1.035692936477183
$ws.Range("L20").Value = 1.040054544762124
$ws.Range("M20").Value = 1.049242714715148
$ws.Range("N20").Value = 1.015185181774554

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026398028853237
$ws.Range("D21").Value = 1.031642351600444
$ws.Range("E21").Value = 1.035807081742263
$ws.Range("F21").Value = 1.044748686842524
$ws.Range("I21").Value = 1.033882324700144
$ws.Range("J21").Value = 1.033053719483462
$ws.Range("K21").Value = 1.035238783011141
$ws.Range("L21").Value = 1.039387937979984
$ws.Range("M21").Value = 1.048296567236888
$ws.Range("N21").Value = 1.015019610477925

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025835352021955
$ws.Range("D22").Value = 1.031225308349454
$ws.Range("E22").Value = 1.035257183212432
$ws.Range("F22").Value = 1.044024311560656
$ws.Range("I22").Value = 1.033767478739315
$ws.Range("J22").Value = 1.03274205240842
$ws.Range("K22").Value = 1.034952975201922
$ws.Range("L22").Value = 1.038969175173637
$ws.Range("M22").Value = 1.047702682542507
$ws.Range("N22").Value = 1.014915452870978

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026133575302342
$ws.Range("D23").Value = 1.031446343205699
$ws.Range("E23").Value = 1.035548592706358
$ws.Range("F23").Value = 1.044408160036441
$ws.Range("I23").Value = 1.033828453547103
$ws.Range("J23").Value = 1.032907279966347
$ws.Range("K23").Value = 1.035104511895571
$ws.Range("L23").Value = 1.039191134745098
$ws.Range("M23").Value = 1.048017417151598
$ws.Range("N23").Value = 1.014970673986479

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027309025806327
$ws.Range("D24").Value = 1.032317583532796
$ws.Range("E24").Value = 1.036698099399034
$ws.Range("F24").Value = 1.045922798192956
$ws.Range("I24").Value = 1.03406644410012
$ws.Range("J24").Value = 1.033557604630552
$ws.Range("K24").Value = 1.035700545862571
$ws.Range("L24").Value = 1.040065726504265
$ws.Range("M24").Value = 1.049258593769846
$ws.Range("N24").Value = 1.015187956611752

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02867585364358
$ws.Range("D25").Value = 1.033330711607091
$ws.Range("E25").Value = 1.038036643819079
$ws.Range("F25").Value = 1.047687561498542
$ws.Range("I25").Value = 1.034338324037567
$ws.Range("J25").Value = 1.034311893112429
$ws.Range("K25").Value = 1.036391021437099
$ws.Range("L25").Value = 1.041082160517494
$ws.Range("M25").Value = 1.050703194293976
$ws.Range("N25").Value = 1.015439840984556

